$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as Text,
# matching the source data (prices/strings stored as text in the sheet),
# without leaving a lasting number-format override on the cell.
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

$ws.Range("D2").Value = '66.368.89'
$ws.Range("D3").Value = '3.522.45'
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("E4").Value = '  +0.05%  '
Set-TextValue "D5" '608.10'
$ws.Range("E5").Value = '  +0.86%  '
Set-TextValue "D6" '145.21'
$ws.Range("E6").Value = '  -1.57%  '
$ws.Range("D7").Value = '3.522.96'
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -0.41%  '
$ws.Range("E10").Value = '  -4.11%  '
Set-TextValue "D11" '8.02'
$ws.Range("E11").Value = '  +2.69%  '
$ws.Range("E12").Value = '  -2.25%  '
$ws.Range("D13").Value = '4.109.14'
$ws.Range("E13").Value = '  +0.17%  '
Set-TextValue "D14" '0.0000207'
$ws.Range("E14").Value = '  -3.42%  '
Set-TextValue "D15" '30.38'
$ws.Range("E15").Value = '  -3.67%  '
$ws.Range("D16").Value = '3.516.76'
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("D17").Value = '66.369.98'
$ws.Range("E17").Value = '  -0.72%  '
$ws.Range("E18").Value = '  -0.03%  '
Set-TextValue "D19" '10.70'
$ws.Range("E19").Value = '  +0.60%  '
Set-TextValue "D20" '6.21'
$ws.Range("E20").Value = '  -2.95%  '
Set-TextValue "D21" '14.91'
$ws.Range("E21").Value = '  -3.17%  '
Set-TextValue "D22" '426.20'
$ws.Range("E22").Value = '  -2.18%  '
$ws.Range("E23").Value = '  -1.68%  '
Set-TextValue "D24" '78.16'
$ws.Range("E24").Value = '  -2.13%  '
$ws.Range("D25").Value = '3.650.55'
$ws.Range("E25").Value = '  +0.20%  '
$ws.Range("E26").Value = '  -0.41%  '
$ws.Range("E27").Value = '  -0.33%  '
Set-TextValue "D28" '9.30'
$ws.Range("E28").Value = '  -5.57%  '
Set-TextValue "D29" '8.02'
$ws.Range("E29").Value = '  -3.22%  '
$ws.Range("E30").Value = '  -0.94%  '
Set-TextValue "D31" '1.00'
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("E32").Value = '  -0.52%  '
$ws.Range("E33").Value = '  -7.52%  '
Set-TextValue "D34" '25.28'
$ws.Range("E34").Value = '  -0.28%  '
$ws.Range("D35").Value = '3.498.43'
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("E37").Value = '  -3.24%  '
$ws.Range("E38").Value = '  -3.01%  '
$ws.Range("E39").Value = '  -4.66%  '
$ws.Range("E40").Value = '  -0.04%  '
Set-TextValue "D41" '170.65'
$ws.Range("E41").Value = '  +0.90%  '
Set-TextValue "D42" '0.0859'
$ws.Range("E42").Value = '  -3.48%  '
Set-TextValue "D43" '5.19'
$ws.Range("E43").Value = '  -4.49%  '
Set-TextValue "D44" '0.889'
$ws.Range("E44").Value = '  -1.00%  '
$ws.Range("E45").Value = '  -9.00%  '
Set-TextValue "D46" '45.49'
$ws.Range("E46").Value = '  -0.49%  '
$ws.Range("E47").Value = '  -7.73%  '
Set-TextValue "D48" '25.85'
$ws.Range("E48").Value = '  -10.90%  '
Set-TextValue "D49" '2.43'
$ws.Range("E49").Value = '  +0.72%  '
Set-TextValue "D50" '7.19'
$ws.Range("E50").Value = '  -3.75%  '
Set-TextValue "D51" '0.949'
$ws.Range("E51").Value = '  -3.43%  '
